$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("B2").Value = 3239
$ws.Range("B3").Value = 1967
$ws.Range("B4").Value = 597
